$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40-91 down to 41-92.
$ws.Rows("40:40").Insert()

# The new row 40 duplicates what is now row 41 (the original row 40 data),
# except the Fecha (date, column D) which becomes 2022-03-09 (serial 44629).
$ws.Range("A40:R40").Value2 = $ws.Range("A41:R41").Value2
$ws.Range("D40").Value2 = 44629
